$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-valued Price cells (column D) keep their original text
# representation (e.g. "1.001", "29.833.77") instead of being parsed as
# numbers, by forcing a Text number format before assigning the value.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.833.77"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +8.54%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.951.42"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +6.77%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.29%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "342.43"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.77%  "

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.19%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4773"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +4.53%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4143"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +8.08%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.91"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +3.89%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08246"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +4.96%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +8.16%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.72"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +7.90%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.943.00"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +8.22%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.187"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +5.97%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.416"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +5.20%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "92.26"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +3.11%  "

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.22%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001061"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +3.82%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06695"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.93%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.08"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +5.75%  "

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.24%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "29.791.64"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +8.49%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.584"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +5.64%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.28"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +4.50%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.254"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.34%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.178.77"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +8.28%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "161.63"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.46%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.21"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +4.28%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.184"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +6.91%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.706"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +8.26%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "122.66"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +4.01%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +9.09%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09630"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.88%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.479"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +12.28%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.690"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +3.30%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.526"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +5.80%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06317"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +7.06%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02320"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +6.14%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.513"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +5.35%  "

$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6110"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +6.60%  "

$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.187"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +3.58%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "10.73"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +8.32%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.001"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.18%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1893"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +4.11%  "

$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.377"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +31.50%  "

$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "WEMIXTOKEN"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.258"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.40%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5719"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +5.92%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "12.47"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +5.98%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.990"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +5.54%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07340"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +5.82%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "114.00"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +3.37%  "
